$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the two anchor paragraphs we need, by their current text, so the
# script is robust to whatever paragraph indices the document starts with.
#   - "Fuente:" is immediately followed by a blank (sz=24) paragraph.
#   - That blank paragraph is immediately followed by the "Botones:"
#     paragraph, whose paragraph mark carries a single underline.
# ---------------------------------------------------------------------------
$blankBeforeBotonesIndex = 0
$botonesIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13)
    if ($t -eq "Botones:") {
        $botonesIndex = $i
        $blankBeforeBotonesIndex = $i - 1
    }
}

# ---------------------------------------------------------------------------
# Step 1: insert the five new list paragraphs right after the blank
# paragraph that precedes "Botones:". Cloning that blank paragraph's
# paragraph mark (plain sz=24, no underline) keeps the new paragraphs free
# of the underline that sits on the old "Botones:" paragraph mark.
# ---------------------------------------------------------------------------
$labels = @("Lista de Formularios:", "Login", "Registro Usuario", "Registro Empresario", "Pantalla Principal")
$insertAfterIndex = $blankBeforeBotonesIndex
foreach ($label in $labels) {
    $srcPara = $d.Paragraphs.Item($insertAfterIndex)
    $srcPara.Range.InsertParagraphAfter()
    $insertAfterIndex = $insertAfterIndex + 1
    $newPara = $d.Paragraphs.Item($insertAfterIndex)
    $newPara.Range.Text = $label
}

# ---------------------------------------------------------------------------
# Step 2: remove the old "Botones:" paragraph (text + underlined paragraph
# mark). Its index shifted down by the number of paragraphs just inserted.
# ---------------------------------------------------------------------------
$oldBotonesIndex = $insertAfterIndex + 1
$d.Paragraphs.Item($oldBotonesIndex).Range.Delete()

# ---------------------------------------------------------------------------
# Step 3: remove the blank sz=32 paragraph that used to sit between
# "Botones:" and the bookmark paragraph - it doesn't survive in the new
# layout (the bookmark paragraph becomes sz=24 and two fresh sz=32 blank
# paragraphs are appended after it instead).
# ---------------------------------------------------------------------------
$blankSz32Index = $oldBotonesIndex
$d.Paragraphs.Item($blankSz32Index).Range.Delete()

# ---------------------------------------------------------------------------
# Step 4: append two new blank sz=32 paragraphs after the bookmark
# paragraph, while it is still sz=32 itself, so the clones pick up sz=32
# too (InsertParagraphAfter clones the calling paragraph's current mark
# formatting).
# ---------------------------------------------------------------------------
$bookmarkIndex = $blankSz32Index
$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)
$bookmarkPara.Range.InsertParagraphAfter()
$d.Paragraphs.Item($bookmarkIndex + 1).Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Step 5: the paragraph that holds the _GoBack bookmark changes from sz=32
# to sz=24. Its range is collapsed (paragraph mark only) so direct Font
# writes on it are ignored; work around that by typing a throw-away
# character ahead of the bookmark, sizing it, then deleting the character
# again, which leaves the size change stamped on the paragraph mark.
# ---------------------------------------------------------------------------
$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)
$insPoint = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start)
$insPoint.InsertBefore("x")
$bookmarkPara = $d.Paragraphs.Item($bookmarkIndex)
$bookmarkPara.Range.Font.Size = 12
$bookmarkPara.Range.Font.SizeBi = 12
$placeholder = $d.Range($bookmarkPara.Range.Start, $bookmarkPara.Range.Start + 1)
$placeholder.Text = ""

Write-Host "Final paragraph dump:"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Host $i ":" "[$($p.Range.Text)]" "size=" $p.Range.Font.Size
}
